$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"
$rng.Value = "2012-06-17"
$rng.ClearFormats()
